$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (A357) down through the new rows (A358:A366)
$ws.Range("A357").Copy()
$ws.Range("A358:A366").PasteSpecial(-4122)

# Fill in the new data rows (dates 2021-08-24 through 2021-09-01)
$ws.Range("A358").Value = 44432
$ws.Range("B358").Value = 5
$ws.Range("C358").Value = 17
$ws.Range("D358").Value = 70.47508498466131

$ws.Range("A359").Value = 44433
$ws.Range("B359").Value = 0
$ws.Range("C359").Value = 17
$ws.Range("D359").Value = 70.47508498466131

$ws.Range("A360").Value = 44434
$ws.Range("B360").Value = 1
$ws.Range("C360").Value = 17
$ws.Range("D360").Value = 70.47508498466131

$ws.Range("A361").Value = 44435
$ws.Range("B361").Value = 3
$ws.Range("C361").Value = 16
$ws.Range("D361").Value = 66.32949175026947

$ws.Range("A362").Value = 44436
$ws.Range("B362").Value = 0
$ws.Range("C362").Value = 16
$ws.Range("D362").Value = 66.32949175026947

$ws.Range("A363").Value = 44437
$ws.Range("B363").Value = 5
$ws.Range("C363").Value = 14
$ws.Range("D363").Value = 58.03830528148578

$ws.Range("A364").Value = 44438
$ws.Range("B364").Value = 1
$ws.Range("C364").Value = 15
$ws.Range("D364").Value = 62.18389851587763

$ws.Range("A365").Value = 44439
$ws.Range("B365").Value = 4
$ws.Range("C365").Value = 14
$ws.Range("D365").Value = 58.03830528148578

$ws.Range("A366").Value = 44440
$ws.Range("B366").Value = 0
$ws.Range("C366").Value = 14
$ws.Range("D366").Value = 58.03830528148578

